$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final JN - Ano / Total / TotalPonto data (years 2011-2020), replacing the
# previous 2015-2020 rows and extending the table through row 11.
$data = @(
  @("2011","21.938","21.938"),
  @("2012","63.615","63.615"),
  @("2013","22.169","22.169"),
  @("2014","58.877","58.877"),
  @("2015","11.804","11.804"),
  @("2016","46.031","46.031"),
  @("2017","27.422","27.422"),
  @("2018","73.156","73.156"),
  @("2019","6.968","6.968"),
  @("2020","42.495","42.495")
)

# Give the new rows (8-11) the same look (fill/border/row height) as the
# existing data rows before writing their values.
$ws.Range("A7:C7").Copy() | Out-Null
$ws.Range("A8:C11").PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(8).RowHeight = 13.073974609375
$ws.Rows.Item(9).RowHeight = 13.073974609375
$ws.Rows.Item(10).RowHeight = 13.073974609375
$ws.Rows.Item(11).RowHeight = 13.073974609375

# Store the data as plain text so that year numbers and dot-separated
# figures (e.g. "21.938") are preserved exactly, instead of being
# reinterpreted as numbers/dates by Excel.
$ws.Range("A2:C11").NumberFormat = "@"

for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $i + 2
  $ws.Cells.Item($r, 1).Value = $data[$i][0]
  $ws.Cells.Item($r, 2).Value = $data[$i][1]
  $ws.Cells.Item($r, 3).Value = $data[$i][2]
}
